$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new data row (row 34): date 2019-11-13 -> serial 43782, kilométrage 1170
# Copy row 33's formatting down first so the new date cell picks up the
# same custom date style, then write the values.
$ws.Range("A33").Copy()
$ws.Range("A34").PasteSpecial(-4122)

$ws.Range("A34").Value = 43782
$ws.Range("B34").Value = 1170

# Update view state to match: scrolled down one row, selection moved to B35
$excel.ActiveWindow.ScrollRow = 14
$ws.Range("B35").Select()
